$d = $word.ActiveDocument

# Find the "Preparation" heading paragraph.
$prepIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Preparation") {
        $prepIndex = $i
        break
    }
}

if ($prepIndex -eq -1) {
    throw "Could not find the 'Preparation' paragraph"
}

# Find the "Solutions" heading paragraph that follows it - the three
# paragraphs in between (the two numbered prep questions and the
# stray "-->" paragraph) are the ones being removed.
$solIndex = -1
for ($i = $prepIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Solutions") {
        $solIndex = $i
        break
    }
}

if ($solIndex -eq -1) {
    throw "Could not find the 'Solutions' paragraph after 'Preparation'"
}

if ($solIndex -gt $prepIndex + 1) {
    $startPara = $d.Paragraphs($prepIndex + 1)
    $endPara = $d.Paragraphs($solIndex - 1)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
